$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 74 updates
$ws.Range("B74").Value = 22232
$ws.Range("C74").Value = 13109
$ws.Range("D74").Value = 12090
$ws.Range("E74").Value = 4777
$ws.Range("F74").Value = 6846
$ws.Range("M74").Value = 2905
$ws.Range("N74").Value = 2755
$ws.Range("O74").Value = 489
$ws.Range("P74").Value = 53
$ws.Range("Q74").Value = 45
$ws.Range("R74").Value = 485
$ws.Range("T74").Value = 143
$ws.Range("U74").Value = 1350
$ws.Range("V74").Value = 30
$ws.Range("W74").Value = 109
$ws.Range("X74").Value = 9
$ws.Range("AB74").Value = 6218
$ws.Range("AC74").Value = 2439
$ws.Range("AF74").Value = 1178
$ws.Range("AG74").Value = 107
$ws.Range("AT74").Value = 89
$ws.Range("AU74").Value = 520
$ws.Range("AV74").Value = 169
$ws.Range("AW74").Value = 57
$ws.Range("BA74").Value = 653
$ws.Range("BF74").Value = 1226
$ws.Range("BH74").Value = 144
$ws.Range("BK74").Value = 259
$ws.Range("BM74").Value = 192
$ws.Range("BN74").Value = 48

# Row 75 updates
$ws.Range("B75").Value = 23234
$ws.Range("C75").Value = 14899
$ws.Range("D75").Value = 13636
$ws.Range("E75").Value = 5122
$ws.Range("F75").Value = 7869
$ws.Range("J75").Value = 75
$ws.Range("O75").Value = 413
$ws.Range("V75").Value = 2
$ws.Range("X75").Value = 59
$ws.Range("AB75").Value = 6935
$ws.Range("AC75").Value = 2643
$ws.Range("AG75").Value = 92
$ws.Range("AO75").Value = 113
$ws.Range("AU75").Value = 616
$ws.Range("AW75").Value = 82
$ws.Range("BA75").Value = 785
$ws.Range("BF75").Value = 1405
$ws.Range("BH75").Value = 181
$ws.Range("BK75").Value = 332
$ws.Range("BM75").Value = 246
$ws.Range("BN75").Value = 61
$ws.Range("BP75").Value = 444
$ws.Range("BQ75").Value = 84
$ws.Range("BR75").Value = 208
$ws.Range("BS75").Value = 151
$ws.Range("BT75").Value = 238
